# Updates cryptos.xlsx price/volume figures (GitHub Actions scheduled refresh).
# Mirrors the upstream diff: per-row Price (D) / Volume(1h) (E) changes, plus
# three pairs of rows whose ranking order swapped (Toncoin/PancakeSwap,
# FTXToken/FraxShare) which carries Coin name + Link along with them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few updated Price values (e.g. "59.60") round-trip as numbers and would
# silently drop a trailing zero (59.60 -> 59.6) if Excel auto-detects them as
# numeric on assignment, so those specific cells are forced to Text format
# first, matching the original text-formatted cell content exactly.
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "37.095.33"
$ws.Cells.Item(2, 5).Value = "  -1.38%  "
$ws.Cells.Item(3, 4).Value = "1.984.28"
$ws.Cells.Item(3, 5).Value = "  -2.77%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.13%  "
$ws.Cells.Item(5, 4).Value = "243.84"
$ws.Cells.Item(5, 5).Value = "  -5.25%  "
$ws.Cells.Item(6, 4).Value = "0.601"
$ws.Cells.Item(6, 5).Value = "  -3.69%  "
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
$ws.Cells.Item(8, 4).Value = "54.56"
$ws.Cells.Item(8, 5).Value = "  -5.09%  "
$ws.Cells.Item(9, 4).Value = "59.60"
$ws.Cells.Item(9, 5).Value = "  +4.14%  "
$ws.Cells.Item(10, 4).Value = "0.371"
$ws.Cells.Item(10, 5).Value = "  -4.25%  "
$ws.Cells.Item(11, 4).Value = "0.0752"
$ws.Cells.Item(11, 5).Value = "  -5.79%  "
$ws.Cells.Item(12, 4).Value = "0.0984"
$ws.Cells.Item(12, 5).Value = "  -4.54%  "
$ws.Cells.Item(13, 4).Value = "2.275.23"
$ws.Cells.Item(13, 5).Value = "  -2.74%  "
$ws.Cells.Item(14, 4).Value = "13.97"
$ws.Cells.Item(14, 5).Value = "  -5.65%  "
$ws.Cells.Item(15, 4).Value = "20.94"
$ws.Cells.Item(15, 5).Value = "  -1.86%  "
$ws.Cells.Item(16, 4).Value = "0.753"
$ws.Cells.Item(16, 5).Value = "  -8.29%  "
$ws.Cells.Item(17, 4).Value = "5.05"
$ws.Cells.Item(17, 5).Value = "  -6.20%  "
$ws.Cells.Item(18, 4).Value = "1.999.30"
$ws.Cells.Item(18, 5).Value = "  -1.96%  "
$ws.Cells.Item(19, 4).Value = "36.999.13"
$ws.Cells.Item(19, 5).Value = "  -1.22%  "
$ws.Cells.Item(20, 4).Value = "68.29"
$ws.Cells.Item(20, 5).Value = "  -2.64%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0807"
$ws.Cells.Item(21, 5).Value = "  -5.69%  "
$ws.Cells.Item(22, 4).Value = "228.54"
$ws.Cells.Item(22, 5).Value = "  -0.27%  "
$ws.Cells.Item(23, 4).Value = "4.96"
$ws.Cells.Item(23, 5).Value = "  -5.25%  "
$ws.Cells.Item(24, 5).Value = "  +0.04%  "
$ws.Cells.Item(25, 2).Value = "PancakeSwap"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(25, 4).Value = "2.41"
$ws.Cells.Item(25, 5).Value = "  -9.65%  "
$ws.Cells.Item(26, 2).Value = "Toncoin"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(26, 4).Value = "2.34"
$ws.Cells.Item(26, 5).Value = "  -0.18%  "
$ws.Cells.Item(27, 4).Value = "161.08"
$ws.Cells.Item(27, 5).Value = "  -1.35%  "
$ws.Cells.Item(28, 4).Value = "8.65"
$ws.Cells.Item(28, 5).Value = "  -5.81%  "
$ws.Cells.Item(29, 4).Value = "19.07"
$ws.Cells.Item(29, 5).Value = "  -4.42%  "
$ws.Cells.Item(30, 4).Value = "0.124"
$ws.Cells.Item(30, 5).Value = "  -10.87%  "
$ws.Cells.Item(31, 4).Value = "1.30"
$ws.Cells.Item(31, 5).Value = "  -4.04%  "
$ws.Cells.Item(32, 5).Value = "  -3.10%  "
$ws.Cells.Item(33, 4).Value = "4.40"
$ws.Cells.Item(33, 5).Value = "  -7.46%  "
$ws.Cells.Item(34, 4).Value = "0.0615"
$ws.Cells.Item(34, 5).Value = "  -7.69%  "
$ws.Cells.Item(35, 4).Value = "4.23"
$ws.Cells.Item(35, 5).Value = "  -6.45%  "
$ws.Cells.Item(36, 4).Value = "2.33"
$ws.Cells.Item(36, 5).Value = "  -7.11%  "
$ws.Cells.Item(37, 4).Value = "0.999"
$ws.Cells.Item(37, 5).Value = "  -0.06%  "
$ws.Cells.Item(38, 4).Value = "1.79"
$ws.Cells.Item(38, 5).Value = "  -1.56%  "
$ws.Cells.Item(39, 4).Value = "3.32"
$ws.Cells.Item(39, 5).Value = "  -4.36%  "
$ws.Cells.Item(40, 4).Value = "5.21"
$ws.Cells.Item(40, 5).Value = "  -3.02%  "
$ws.Cells.Item(41, 5).Value = "  -0.32%  "
$ws.Cells.Item(42, 4).Value = "1.421.00"
$ws.Cells.Item(42, 5).Value = "  +0.78%  "
$ws.Cells.Item(43, 4).Value = "1.13"
$ws.Cells.Item(43, 5).Value = "  -5.13%  "
$ws.Cells.Item(44, 4).Value = "0.0203"
$ws.Cells.Item(44, 5).Value = "  -6.47%  "
$ws.Cells.Item(45, 4).Value = "0.0888"
$ws.Cells.Item(45, 5).Value = "  -8.15%  "
$ws.Cells.Item(46, 4).Value = "87.77"
$ws.Cells.Item(46, 5).Value = "  -4.08%  "
$ws.Cells.Item(47, 4).Value = "15.31"
$ws.Cells.Item(47, 5).Value = "  -5.53%  "
$ws.Cells.Item(48, 4).Value = "0.998"
$ws.Cells.Item(48, 5).Value = "  -5.21%  "
$ws.Cells.Item(49, 4).Value = "2.86"
$ws.Cells.Item(49, 5).Value = "  -0.71%  "
$ws.Cells.Item(50, 2).Value = "FraxShare"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(50, 4).Value = "6.63"
$ws.Cells.Item(50, 5).Value = "  -11.00%  "
$ws.Cells.Item(51, 2).Value = "FTXToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(51, 4).Value = "3.58"
$ws.Cells.Item(51, 5).Value = "  +10.86%  "

Write-Host "Applied updates to cryptos sheet"
